$wb = $excel.ActiveWorkbook

# --- doordash sheet (sheet1) ---
$ws1 = $wb.Worksheets.Item("doordash")
$ws1.Range("C4").Value = 68244709.11999999
$ws1.Range("C5").Value = 64339464.20000002
$ws1.Range("C6").Value = 72433065.29000002
$ws1.Range("C8").Value = 73103432.77999997
$ws1.Range("C9").Value = 77135047.39999998
$ws1.Range("C10").Value = 79303082.89999998
$ws1.Range("C11").Value = 85827325.39999998
$ws1.Range("C12").Value = 84346802.20000005
$ws1.Range("C13").Value = 87457261.29999995
$ws1.Range("C15").Value = 95192460.33999999
$ws1.Range("C16").Value = 111952800.2
$ws1.Range("C21").Value = 237625455.6099999
$ws1.Range("C22").Value = 229289794
$ws1.Range("C23").Value = 255090913
$ws1.Range("C24").Value = 248271130
$ws1.Range("C25").Value = 276198421
$ws1.Range("C26").Value = 314656599.8
$ws1.Range("C27").Value = 278654048.6
$ws1.Range("C28").Value = 319326490.2
$ws1.Range("C29").Value = 323398158.4
$ws1.Range("C30").Value = 333916389
$ws1.Range("C31").Value = 306876533.6900003
$ws1.Range("C32").Value = 321985059.3099997
$ws1.Range("C33").Value = 322088494.1600003
$ws1.Range("C34").Value = 313324738.8399997
$ws1.Range("C35").Value = 336797364
$ws1.Range("B36").Value = 9626952.489999995
$ws1.Range("C36").Value = 321331607
$ws1.Range("B37").Value = 10134213.59999999
$ws1.Range("C37").Value = 343800779
$ws1.Range("C38").Value = 365342652.3
$ws1.Range("C39").Value = 330918924.8999999
$ws1.Range("C40").Value = 370481532.8000001
$ws1.Range("C41").Value = 370048782
$ws1.Range("C42").Value = 381151540.6199999
$ws1.Range("C43").Value = 353226014.3800001
$ws1.Range("C44").Value = 364289812
$ws1.Range("C45").Value = 356142624
$ws1.Range("B46").Value = 10308955.07000001
$ws1.Range("C46").Value = 353876016
$ws1.Range("B47").Value = 10869381.41999999
$ws1.Range("C47").Value = 374789342
$ws1.Range("B48").Value = 10425674.10000001
$ws1.Range("C48").Value = 358899237.73
$ws1.Range("B49").Value = 10898187.13
$ws1.Range("C50").Value = 385847187.4
$ws1.Range("C51").Value = 370191995.7
$ws1.Range("C52").Value = 416408720.9
$ws1.Range("C53").Value = 404036293.4299998
$ws1.Range("C54").Value = 411526873.5700002
$ws1.Range("C55").Value = 396539587
$ws1.Range("C56").Value = 407758991
$ws1.Range("C57").Value = 395753325
$ws1.Range("C58").Value = 395633410.2700005
$ws1.Range("B60").Value = 11295149.90000001
$ws1.Range("B61").Value = 12068675.09999999
$ws1.Range("C63").Value = 417453068.72
$ws1.Range("C64").Value = 462026821.4
$ws1.Range("C65").Value = 439918083.8000002
$ws1.Range("C68").Value = 442249556.3899994
$ws1.Range("C69").Value = 447279830
$ws1.Range("B70").Value = 12398684.09999999
$ws1.Range("C70").Value = 432660425
$ws1.Range("B71").Value = 12914607.60000001
$ws1.Range("C71").Value = 450245500
$ws1.Range("B72").Value = 12747305
$ws1.Range("C72").Value = 449028645
$ws1.Range("B73").Value = 13184010.79999998
$ws1.Range("C73").Value = 473612484.8899994
$ws1.Range("C74").Value = 473375846.9
$ws1.Range("C75").Value = 451204963.2900001
$ws1.Range("B81").Value = 14325866.47
$ws1.Range("C81").Value = 516968612.2300005

# --- ubereats sheet (sheet2) ---
$ws2 = $wb.Worksheets.Item("ubereats")
# Add new column D header ("Unnamed: 3") matching style of C1
$ws2.Range("C1").Copy($ws2.Range("D1"))
$ws2.Range("D1").Value = "Unnamed: 3"

$ws2.Range("C4").Value = 54012045.04000001
$ws2.Range("C5").Value = 52282240.45999998
$ws2.Range("C7").Value = 39943646.09000003
$ws2.Range("C8").Value = 31431968.29999995
$ws2.Range("C9").Value = 39897479
$ws2.Range("C10").Value = 51902411.10000002
$ws2.Range("C11").Value = 54407014.89999998
$ws2.Range("C12").Value = 54672489.20000005
$ws2.Range("C13").Value = 57042329.54000008
$ws2.Range("C20").Value = 122836415.51
$ws2.Range("C21").Value = 127960388
$ws2.Range("C23").Value = 136129154.3599999
$ws2.Range("C24").Value = 133846309.5900002
$ws2.Range("C25").Value = 135404066.4099998
$ws2.Range("C28").Value = 156113589.2899999
$ws2.Range("C29").Value = 152682319.17
$ws2.Range("C30").Value = 159059851.73
$ws2.Range("C31").Value = 156111480.83
$ws2.Range("C33").Value = 169907582.61
$ws2.Range("C34").Value = 163862509.8400002
$ws2.Range("C35").Value = 177597265.1599998
$ws2.Range("C36").Value = 166229189.1099999
$ws2.Range("C38").Value = 186130025.4
$ws2.Range("C40").Value = 188891670.6
$ws2.Range("C41").Value = 184212486.1500001
$ws2.Range("C43").Value = 175875548.2299998
$ws2.Range("C44").Value = 183771368
$ws2.Range("C45").Value = 179451891
$ws2.Range("B46").Value = 6019439.979999997
$ws2.Range("C46").Value = 176514471
$ws2.Range("B47").Value = 6331139.729999997
$ws2.Range("C47").Value = 185736573
$ws2.Range("B48").Value = 5578565.350000009
$ws2.Range("C48").Value = 174064280
$ws2.Range("B49").Value = 5411334.799999997
$ws2.Range("C49").Value = 182863162
$ws2.Range("C52").Value = 198821285.7699999
$ws2.Range("C53").Value = 193013656.7
$ws2.Range("C54").Value = 194812830.5700002
$ws2.Range("C55").Value = 189554022.4299998
$ws2.Range("C56").Value = 194728052
$ws2.Range("C57").Value = 190577162.8900001
$ws2.Range("C58").Value = 195160018.1099999
$ws2.Range("C59").Value = 202385409
$ws2.Range("C60").Value = 187986616
$ws2.Range("C61").Value = 199707635
$ws2.Range("C63").Value = 198220762.89
$ws2.Range("C65").Value = 209985530.8500001
$ws2.Range("C67").Value = 215331303.65
$ws2.Range("C68").Value = 218760042.9500003
$ws2.Range("C69").Value = 238653540.0499997
$ws2.Range("C70").Value = 237419822.5999999
$ws2.Range("C71").Value = 240657337.4000001
$ws2.Range("C72").Value = 237233258
$ws2.Range("C73").Value = 247731498.6800003
$ws2.Range("C76").Value = 266909264.1299999
$ws2.Range("C77").Value = 252389623.11
$ws2.Range("C78").Value = 270334393.09
$ws2.Range("C79").Value = 257429312
$ws2.Range("C80").Value = 258224971.4600005
$ws2.Range("C81").Value = 267283138.5399995
